$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 183; this shifts the existing rows
# 183:255 down to 184:256 (and keeps column D's date number format on
# the newly inserted row, same as the rows around it).
$ws.Rows(183).Insert()

# Populate the newly inserted row 183 with the new record's data.
# Columns A, B, C, E, F, G, H, I, J, K are constant across this block
# of data (same Mercado / Region / Producto), so copy them from the
# row directly below (the old row 183, now shifted to row 184).
$ws.Range("A183").Value = $ws.Range("A184").Value()
$ws.Range("B183").Value = $ws.Range("B184").Value()
$ws.Range("C183").Value = $ws.Range("C184").Value()
$ws.Range("D183").Value = 44917
$ws.Range("E183").Value = $ws.Range("E184").Value()
$ws.Range("F183").Value = $ws.Range("F184").Value()
$ws.Range("G183").Value = $ws.Range("G184").Value()
$ws.Range("H183").Value = $ws.Range("H184").Value()
$ws.Range("I183").Value = $ws.Range("I184").Value()
$ws.Range("J183").Value = $ws.Range("J184").Value()
$ws.Range("K183").Value = $ws.Range("K184").Value()
$ws.Range("L183").Value = 'Primera'
$ws.Range("M183").Value = 530
$ws.Range("N183").Value = 3000
$ws.Range("O183").Value = 3500
$ws.Range("P183").Value = 3264
$ws.Range("Q183").Value = '$/bandeja 2 kilos'
$ws.Range("R183").Value = 'Región Metropolitana'
$ws.Range("S183").Value = 1632
$ws.Range("T183").Value = 2
